$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Driver Summary")

$ws.Range("B12").Value = 1074341
$ws.Range("C12").Value = 4452
$ws.Range("E12").Value = 1647
$ws.Range("F12").Value = 1078973
